$wb = $excel.ActiveWorkbook

# Select D17 on the Italy sheet to mirror the final recorded selection state.
$italy = $wb.Worksheets.Item("Italy")
$italy.Activate()
$italy.Range("D17").Select()

# Duplicate the Italy sheet, place it right after Italy, and rename it to Spain.
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item($italy.Index + 1)
$spain.Name = "Spain"

# Localize the market name and part number on the new sheet.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2034"

# Make Spain the active sheet/tab with B9 selected, matching the final UI state.
$spain.Activate()
$spain.Range("B9").Select()
